$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the rich-text title in A1: append " - SANTA ROSA" after CACERES ---
#     "...IEP N°54411 "ANDRES AVELINO CACERES",DISTRITO..."
#  -> "...IEP N°54411 "ANDRES AVELINO CACERES - SANTA ROSA",DISTRITO..."
$cell = $ws.Range("A1")
$cell.Characters(129, 7).Text = "CACERES - SANTA ROSA"

# Re-assert the run formatting so the engine keeps the rich-text runs instead
# of collapsing the whole string to a single plain run.
$fMid = $cell.Characters(33, 169).Font
$fMid.Bold = $true
$fMid.Name = "Arial"
$fMid.Size = 10

$fEnd = $cell.Characters(202, 1).Font
$fEnd.Bold = $true
$fEnd.Name = "Arial Narrow"
$fEnd.Size = 10

# --- 2. C11 / D11 change from static half-values to live formulas on F11 ---
$ws.Range("C11").Formula = "=F11/2"
$ws.Range("D11").Formula = "=F11/2"

# --- 3. Selection moves from C9 to the header row A1:F1 ---
$ws.Range("A1:F1").Select() | Out-Null
